# Updates the cryptos price/volume table with the latest scraped values.
# D column cells that look like plain numbers (e.g. "243.68", "1.00") must
# be forced to text via NumberFormat "@" before assignment, otherwise Excel
# auto-converts them to numeric values and silently drops formatting such
# as trailing zeros. The style is reset back to "Normal" right after so the
# cell keeps the workbook's original (unstyled) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.416.27"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.942.31"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.06"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.53%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0831"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("E13").Value = "  -5.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.33%  "
$ws.Range("D15").Value = "2.226.05"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.84%  "
$ws.Range("E17").Value = "  -4.21%  "
$ws.Range("D18").Value = "1.943.85"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").Value = "36.394.56"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.20%  "
$ws.Range("E30").Value = "  -10.05%  "
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  -5.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0631"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  -3.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.32%  "
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.39%  "
$ws.Range("E43").Value = "  -5.30%  "
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.38%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.347.49"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("E48").Value = "  -6.34%  "
$ws.Range("E49").Value = "  -4.89%  "
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.91%  "
